$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.207376956939697
$ws.Range("B1").Value = 3.845653057098389
$ws.Range("C1").Value = 4.05192756652832
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 3.510376453399658
